# Update cryptos list data for Fri Aug 23 07:34:11 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.925.72"
$ws.Range("D3").Value = "2.672.55"
$ws.Range("E3").Value = "  +2.60%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.20"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.15"
$ws.Range("E6").Value = "  +1.77%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("E10").Value = "  +1.79%  "
$ws.Range("E11").Value = "  +4.45%  "
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").Value = "3.144.95"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.74"
$ws.Range("E14").Value = "  +10.87%  "
$ws.Range("D15").Value = "60.916.95"
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("D17").Value = "2.673.46"
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.60"
$ws.Range("E18").Value = "  +2.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.74"
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "351.39"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("E23").Value = "  +1.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.07"
$ws.Range("E24").Value = "  +1.34%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  +1.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.14"
$ws.Range("E27").Value = "  +5.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.98"
$ws.Range("E28").Value = "  +7.40%  "
$ws.Range("D29").Value = "0.0₃0815"
$ws.Range("E29").Value = "  +3.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.80"
$ws.Range("E30").Value = "  +6.71%  "
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.50"
$ws.Range("E32").Value = "  +3.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.89"
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.06"
$ws.Range("E34").Value = "  +8.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.45"
$ws.Range("E35").Value = "  +5.46%  "
$ws.Range("E36").Value = "  +8.51%  "
$ws.Range("E37").Value = "  +3.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "329.61"
$ws.Range("E38").Value = "  +11.89%  "
$ws.Range("E39").Value = "  +4.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.37"
$ws.Range("E40").Value = "  +1.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.880"
$ws.Range("E41").Value = "  +4.29%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.55"
$ws.Range("E42").Value = "  +3.94%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.21"
$ws.Range("E43").Value = "  +5.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "134.61"
$ws.Range("E44").Value = "  -1.91%  "
$ws.Range("E45").Value = "  +1.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0561"
$ws.Range("E47").Value = "  +0.91%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0248"
$ws.Range("E48").Value = "  +3.20%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.46"
$ws.Range("E50").Value = "  +3.59%  "
$ws.Range("D51").Value = "2.119.16"
$ws.Range("E51").Value = "  +4.57%  "
